$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "BLOCO DE CONCRETO DE VEDAÇÃO - BLOCO INTEIRO 14 X 19 X 19 - RESISTÊNICA 3,0 MPA (ALTURA: 190 MM / COMPRIMENTO: 390 MM / LARGURA: 190 MM)"
$ws.Range("A4").Value = "BLOCO DE CONCRETO PARA PAVIMENTAÇÃO INTERTRAVADA RETANGULAR - TRÁFEGO PESADO (COMPRIMENTO: 100 MM / ESPESSURA: 60 MM / LARGURA: 200 MM)"
$ws.Range("A5").Value = "PEDRA BRITADA 2"
$ws.Range("A6").Value = "CIMENTO PORTLAND CP II-Z-32 (RESISTÊNCIA: 32,00 MPA)"

[void]$ws.Range("A6").Select()
